# Reorders the run-per-line entries of the "Requisitos" bullet list
# paragraph so they appear in the order shown in the target revision,
# while keeping each entry as its own run (text + manual line break)
# just like the original markup.
$d = $word.ActiveDocument

# Locate the ListBullet paragraph that follows the "Requisitos" heading.
$targetPara = $null
foreach ($para in $d.Paragraphs) {
    if ($para.Range.Text.TrimEnd([char]13) -eq "Requisitos") {
        $targetPara = $para.Next()
        break
    }
}
if ($targetPara -eq $null) {
    $targetPara = $d.Paragraphs.Last
}

$start = $targetPara.Range.Start
# Length of the existing run content, excluding the trailing paragraph mark.
$oldLen = $targetPara.Range.Text.Length - 1

$newOrder = @(
    "LOB1045 -  Leitura e Produção de Textos Acadêmicos  (Requisito)",
    "LOB1056 -  Introdução aos Métodos Numéricos e Computacionais  (Requisito)",
    "LOQ4095 -  Química Geral Experimental  (Requisito)",
    "LOQ4098 -  Fundamentos de Química para Engenharia II (Requisito)",
    "LOQ4100 -  Fundamentos de Química para Engenharia I (Requisito)",
    "LOB1006 -  Cálculo IV  (Requisito)",
    "LOB1037 -  Àlgebra Linear  (Requisito)",
    "LOB1040 -  Laboratório de Eletricidade  (Requisito)",
    "LOB1053 -  Física III  (Requisito)",
    "LOB1003 -  Cálculo I  (Requisito)",
    "LOB1009 -  Leitura e Interpretação de Desenho Técnico  (Requisito)",
    "LOB1011 -  Eletricidade Aplicada  (Requisito)",
    "LOB1018 -  Física I  (Requisito)",
    "LOB1024 -  Mecânica  (Requisito)",
    "LOB1036 -  Geometria Analítica  (Requisito)",
    "LOB1038 -  Física Experimental I  (Requisito)",
    "LOB1039 -  Física Experimental III  (Requisito)",
    "LOB1052 -  Cálculo III  (Requisito)",
    "LOM3016 -  Introdução à  Ciência dos Materiais  (Requisito)",
    "LOB1004 -  Cálculo II  (Requisito)",
    "LOB1019 -  Física II  (Requisito)"
)

# Insert the lines in reverse order, always at the paragraph's original
# start position: each InsertBefore call lands ahead of the previously
# inserted text, so the final on-disk order matches $newOrder. Each call
# creates its own run (the runtime only coalesces runs that are edited in
# place, not freshly inserted text), mirroring the one-run-per-line shape
# of the original document.
for ($i = $newOrder.Length - 1; $i -ge 0; $i--) {
    $insertPoint = $d.Range($start, $start)
    $insertPoint.InsertBefore($newOrder[$i] + [char]11)
}

# The untouched original 21 runs were pushed right after the newly
# inserted text; delete that old block now.
$newLen = 0
foreach ($item in $newOrder) { $newLen = $newLen + $item.Length + 1 }
$oldStart = $start + $newLen
$oldRange = $d.Range($oldStart, $oldStart + $oldLen)
$oldRange.Delete()

Write-Host "Reordered Requisitos list."
